$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2539184952978056
$ws.Range("C2").Value = 0.4263322884012539
$ws.Range("J2").Value = 0.02821316614420063
$ws.Range("P2").Value = 0.1567398119122257
$ws.Range("S2").Value = 0.1347962382445141
$ws.Range("C3").Value = 0.02898550724637681
$ws.Range("J3").Value = 0.03623188405797102
$ws.Range("P3").Value = 0.7246376811594203
$ws.Range("S3").Value = 0.2101449275362319
$ws.Range("J4").Value = 0.1071428571428571
$ws.Range("P4").Value = 0.5714285714285714
$ws.Range("S4").Value = 0.3214285714285715
$ws.Range("B6").Value = 0.04166666666666666
$ws.Range("D6").Value = 0.01893939393939394
$ws.Range("F6").Value = 0.08333333333333333
$ws.Range("J6").Value = 0.2765151515151515
$ws.Range("O6").Value = 0.01136363636363636
$ws.Range("Q6").Value = 0.1553030303030303
$ws.Range("R6").Value = 0.07954545454545454
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.1210762331838565
$ws.Range("D7").Value = 0.02690582959641256
$ws.Range("E7").Value = 0.008968609865470852
$ws.Range("F7").Value = 0.05829596412556054
$ws.Range("J7").Value = 0.1255605381165919
$ws.Range("O7").Value = 0.02242152466367713
$ws.Range("Q7").Value = 0.1838565022421525
$ws.Range("R7").Value = 0.1165919282511211
$ws.Range("S7").Value = 0.336322869955157
$ws.Range("B8").Value = 0.09482758620689655
$ws.Range("D8").Value = 0.01724137931034483
$ws.Range("F8").Value = 0.06465517241379311
$ws.Range("J8").Value = 0.103448275862069
$ws.Range("O8").Value = 0.008620689655172414
$ws.Range("Q8").Value = 0.1745689655172414
$ws.Range("R8").Value = 0.1336206896551724
$ws.Range("S8").Value = 0.4030172413793103
$ws.Range("B9").Value = 0.08267716535433071
$ws.Range("D9").Value = 0.01968503937007874
$ws.Range("E9").Value = 0.007874015748031496
$ws.Range("F9").Value = 0.07874015748031496
$ws.Range("J9").Value = 0.08661417322834646
$ws.Range("O9").Value = 0.03149606299212598
$ws.Range("Q9").Value = 0.1889763779527559
$ws.Range("R9").Value = 0.1102362204724409
$ws.Range("S9").Value = 0.3937007874015748
$ws.Range("B10").Value = 0.09337134711332858
$ws.Range("D10").Value = 0.02280826799714897
$ws.Range("F10").Value = 0.0691375623663578
$ws.Range("J10").Value = 0.1233071988595866
$ws.Range("O10").Value = 0.01069137562366358
$ws.Range("Q10").Value = 0.1967213114754098
$ws.Range("R10").Value = 0.1083392729864576
$ws.Range("S10").Value = 0.3756236635780471
$ws.Range("G11").Value = 0.1257668711656442
$ws.Range("J11").Value = 0.06134969325153374
$ws.Range("K11").Value = 0.1748466257668712
$ws.Range("L11").Value = 0.6165644171779141
$ws.Range("S11").Value = 0.02147239263803681
$ws.Range("G12").Value = 0.7428571428571429
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.009523809523809525
$ws.Range("L12").Value = 0.0380952380952381
$ws.Range("S12").Value = 0.009523809523809525
$ws.Range("G13").Value = 0.6530612244897959
$ws.Range("J13").Value = 0.3265306122448979
$ws.Range("S13").Value = 0.02040816326530612
$ws.Range("G14").Value = 0.5714285714285714
$ws.Range("J14").Value = 0.2857142857142857
$ws.Range("S14").Value = 0.1428571428571428
$ws.Range("F15").Value = 0.01724137931034483
$ws.Range("H15").Value = 0.1810344827586207
$ws.Range("I15").Value = 0.09482758620689655
$ws.Range("J15").Value = 0.3405172413793103
$ws.Range("K15").Value = 0.0603448275862069
$ws.Range("M15").Value = 0.008620689655172414
$ws.Range("O15").Value = 0.08189655172413793
$ws.Range("S15").Value = 0.2155172413793103
$ws.Range("H16").Value = 0.1942857142857143
$ws.Range("I16").Value = 0.1485714285714286
$ws.Range("J16").Value = 0.3657142857142857
$ws.Range("K16").Value = 0.1028571428571429
$ws.Range("M16").Value = 0.03428571428571429
$ws.Range("N16").Value = 0.005714285714285714
$ws.Range("O16").Value = 0.04571428571428571
$ws.Range("S16").Value = 0.1028571428571429
$ws.Range("F17").Value = 0.02268041237113402
$ws.Range("H17").Value = 0.154639175257732
$ws.Range("I17").Value = 0.1195876288659794
$ws.Range("J17").Value = 0.4391752577319588
$ws.Range("K17").Value = 0.08453608247422681
$ws.Range("M17").Value = 0.02268041237113402
$ws.Range("N17").Value = 0.002061855670103093
$ws.Range("O17").Value = 0.06597938144329897
$ws.Range("S17").Value = 0.088659793814433
$ws.Range("F18").Value = 0.03460207612456748
$ws.Range("H18").Value = 0.1730103806228374
$ws.Range("I18").Value = 0.09688581314878893
$ws.Range("J18").Value = 0.4740484429065744
$ws.Range("K18").Value = 0.07612456747404844
$ws.Range("M18").Value = 0.01384083044982699
$ws.Range("N18").Value = 0.006920415224913495
$ws.Range("O18").Value = 0.03806228373702422
$ws.Range("S18").Value = 0.08650519031141868
$ws.Range("F19").Value = 0.02160953800298063
$ws.Range("H19").Value = 0.1929955290611028
$ws.Range("I19").Value = 0.09016393442622951
$ws.Range("J19").Value = 0.3651266766020864
$ws.Range("K19").Value = 0.1244411326378539
$ws.Range("M19").Value = 0.01937406855439642
$ws.Range("N19").Value = 0.002235469448584203
$ws.Range("O19").Value = 0.07153502235469449
